$wb = $excel.ActiveWorkbook

$wsOrg = $wb.Worksheets.Item("OrgData")
$wsOrg.Range("A2").Value = "AUTO_ORG_PAVNO"

$wsSubOrg = $wb.Worksheets.Item("SubOrgData")
$wsSubOrg.Range("A2").Value = "AUTO_SUB_ORG_ESVEA"
